$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.321
$ws.Range("G2").Value = -69.66666666666667
$ws.Range("H2").Value = -69.66666666666667
$ws.Range("I2").Value = -72.76190476190476
$ws.Range("J2").Value = -72.76190476190476
$ws.Range("K2").Value = -1.842
$ws.Range("L2").Value = -87.71428571428571
$ws.Range("U2").Value = 0.284
$ws.Range("V2").Value = 0.03761589403973509
$ws.Range("X2").Value = 0.09303633707596873
$ws.Range("Z2").Value = -0.01397205588822356
$ws.Range("AB2").Value = 0.07666615179801008
$ws.Range("AD2").Value = 1.71
$ws.Range("AF2").Value = 1.71
$ws.Range("AG2").Value = 1.426
$ws.Range("AH2").Value = 0.1846652267818575
$ws.Range("AI2").Value = -42.74999999999996
$ws.Range("AJ2").Value = 0.1588680926916221
$ws.Range("AK2").Value = -4.401234567901233
$ws.Range("AL2").Value = 0.578
$ws.Range("AM2").Value = 0.578
$ws.Range("AN2").Value = -1.176875430144529
$ws.Range("AO2").Value = -2.643598615916955
$ws.Range("AP2").Value = -0.981417756366139
$ws.Range("AQ2").Value = -2.643598615916955

# Row 3
$ws.Range("D3").Value = -0.321
$ws.Range("G3").Value = -57.14285714285714
$ws.Range("H3").Value = -57.14285714285714
$ws.Range("I3").Value = -59.04761904761904
$ws.Range("J3").Value = -59.04761904761904
$ws.Range("K3").Value = -1.53
$ws.Range("L3").Value = -72.85714285714285
$ws.Range("U3").Value = 0.284
$ws.Range("V3").Value = 0.1339622641509434
$ws.Range("W3").Value = 0.5503597122302158
$ws.Range("X3").Value = 0.1106395911876437
$ws.Range("Y3").Value = 0.4397201210425721
$ws.Range("Z3").Value = -0.01397205588822356
$ws.Range("AA3").Value = 0.8250166333998671
$ws.Range("AB3").Value = 0.0778992206317264
$ws.Range("AC3").Value = 0.7471174127681407
$ws.Range("AD3").Value = 1.71
$ws.Range("AF3").Value = 1.71
$ws.Range("AG3").Value = 1.426
$ws.Range("AH3").Value = 0.4464751958224543
$ws.Range("AI3").Value = -42.74999999999996
$ws.Range("AJ3").Value = 0.4021432600112803
$ws.Range("AK3").Value = -4.401234567901233
$ws.Range("AL3").Value = 0.578
$ws.Range("AM3").Value = 0.578
$ws.Range("AN3").Value = -1.436974789915966
$ws.Range("AO3").Value = -2.145328719723183
$ws.Range("AP3").Value = -1.198319327731092
$ws.Range("AQ3").Value = -2.145328719723183

# Row 4
$ws.Range("K4").Value = -0.312
$ws.Range("X4").Value = 0.07543308296429375
$ws.Range("AB4").Value = 0.07543308296429375

# Deletions in row 4 (cells removed entirely)
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("L4").ClearContents()
